# Update "Report Activity" per-department metrics in Sheet1 (B2:L18).
# Commit: "alligned hour as index for chart ..> startin from 9 + records of
# pasthour update only at 00. ELiminated initialization" — this refreshes the
# snapshot numbers (counts, averages, service levels, etc.) for every
# department row now that the hour index starts at 9 and the "past hour"
# figures only update at minute :00 (no more stale init values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1161
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 459.123
$ws.Range("E2").Value = 674.1279999999999
$ws.Range("F2").Value = 77
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 6298
$ws.Range("I2").Value = 17287
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 136
$ws.Range("L2").Value = 4936
$ws.Range("B3").Value = 257
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 0
$ws.Range("B4").Value = 173
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 67
$ws.Range("F4").Value = 17
$ws.Range("H4").Value = 2537
$ws.Range("I4").Value = 2800
$ws.Range("J4").Value = 0
$ws.Range("B5").Value = 1514
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 151
$ws.Range("E5").Value = 154
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 372
$ws.Range("I5").Value = 21009
$ws.Range("J5").Value = 2
$ws.Range("B6").Value = 1638
$ws.Range("D6").Value = 279
$ws.Range("E6").Value = 285
$ws.Range("F6").Value = 11
$ws.Range("H6").Value = 3188
$ws.Range("I6").Value = 56183
$ws.Range("J6").Value = 8
$ws.Range("B7").Value = 117
$ws.Range("D7").Value = 50
$ws.Range("E7").Value = 50
$ws.Range("I7").Value = 8750
$ws.Range("B8").Value = 1019
$ws.Range("D8").Value = 814
$ws.Range("E8").Value = 934
$ws.Range("F8").Value = 43
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 679
$ws.Range("I8").Value = 12572
$ws.Range("J8").Value = 16
$ws.Range("K8").Value = 71
$ws.Range("L8").Value = 1079
$ws.Range("B9").Value = 256
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 3333
$ws.Range("B10").Value = 1314
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 469
$ws.Range("E10").Value = 510
$ws.Range("F10").Value = 41
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 382
$ws.Range("I10").Value = 16122
$ws.Range("J10").Value = 9
$ws.Range("B11").Value = 1765
$ws.Range("D11").Value = 227
$ws.Range("E11").Value = 257
$ws.Range("F11").Value = 30
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 2596
$ws.Range("I11").Value = 13186
$ws.Range("J11").Value = 3
$ws.Range("B12").Value = 1376
$ws.Range("D12").Value = 64
$ws.Range("E12").Value = 69
$ws.Range("F12").Value = 5
$ws.Range("H12").Value = 758
$ws.Range("I12").Value = 12213
$ws.Range("J12").Value = 3
$ws.Range("B13").Value = 847
$ws.Range("D13").Value = 402
$ws.Range("E13").Value = 406
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 80
$ws.Range("I13").Value = 12211
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("B14").Value = 258
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 1667
$ws.Range("I14").Value = 8000
$ws.Range("J14").Value = 0
$ws.Range("D15").Value = 135
$ws.Range("E15").Value = 143
$ws.Range("F15").Value = 10
$ws.Range("H15").Value = 787
$ws.Range("I15").Value = 20861
$ws.Range("J15").Value = 3
$ws.Range("B16").Value = 177
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = 4286
$ws.Range("B17").Value = 1745
$ws.Range("D17").Value = 31
$ws.Range("E17").Value = 32
$ws.Range("F17").Value = 1
$ws.Range("H17").Value = 1429
$ws.Range("I17").Value = 15614
$ws.Range("J17").Value = 0
$ws.Range("B18").Value = 4345
$ws.Range("D18").Value = 458
$ws.Range("E18").Value = 1230
$ws.Range("F18").Value = 79
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 2196
$ws.Range("I18").Value = 15134
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = 682
$ws.Range("L18").Value = 16032
